$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.894.35'
$ws.Range('E2').Value = '  +4.57%  '
$ws.Range('D3').Value = '3.559.15'
$ws.Range('E3').Value = '  +3.77%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = "'601.72"
$ws.Range('E5').Value = '  +4.16%  '
$ws.Range('D6').Value = "'171.99"
$ws.Range('E6').Value = '  +4.29%  '
$ws.Range('D7').Value = '3.553.30'
$ws.Range('E7').Value = '  +3.84%  '
$ws.Range('E8').Value = '  +2.21%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').Value = "'0.194"
$ws.Range('E10').Value = '  +6.36%  '
$ws.Range('D11').Value = "'7.35"
$ws.Range('E11').Value = '  +9.68%  '
$ws.Range('D12').Value = "'0.587"
$ws.Range('E12').Value = '  +3.77%  '
$ws.Range('D13').Value = "'46.35"
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').Value = "'0.0000277"
$ws.Range('E14').Value = '  +3.29%  '
$ws.Range('D15').Value = '4.132.30'
$ws.Range('E15').Value = '  +3.27%  '
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').Value = "'610.17"
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '70.876.72'
$ws.Range('E18').Value = '  +4.27%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.557.32'
$ws.Range('E19').Value = '  +3.51%  '
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('E21').Value = '  +1.62%  '
$ws.Range('D22').Value = "'0.879"
$ws.Range('E22').Value = '  +1.53%  '
$ws.Range('D23').Value = "'9.30"
$ws.Range('E23').Value = '  -14.27%  '
$ws.Range('D24').Value = "'15.72"
$ws.Range('E24').Value = '  +2.57%  '
$ws.Range('E25').Value = '  +1.80%  '
$ws.Range('D26').Value = "'3.71"
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = "'2.61"
$ws.Range('E28').Value = '  +2.04%  '
$ws.Range('D29').Value = "'33.99"
$ws.Range('E29').Value = '  +6.51%  '
$ws.Range('D30').Value = "'9.07"
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('D31').Value = "'705.85"
$ws.Range('E31').Value = '  +16.72%  '
$ws.Range('E32').Value = '  +1.39%  '
$ws.Range('E33').Value = '  -0.65%  '
$ws.Range('D34').Value = "'7.07"
$ws.Range('E34').Value = '  +4.95%  '
$ws.Range('E35').Value = '  +1.26%  '
$ws.Range('D36').Value = "'3.63"
$ws.Range('E36').Value = '  +6.95%  '
$ws.Range('D37').Value = "'0.101"
$ws.Range('E37').Value = '  +1.07%  '
$ws.Range('D38').Value = "'10.76"
$ws.Range('E38').Value = '  +1.82%  '
$ws.Range('D39').Value = "'0.0479"
$ws.Range('E39').Value = '  +11.13%  '
$ws.Range('D40').Value = "'56.95"
$ws.Range('E40').Value = '  +0.66%  '
$ws.Range('D41').Value = "'0.999"
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('D42').Value = "'0.144"
$ws.Range('E42').Value = '  +7.33%  '
$ws.Range('D43').Value = '3.369.39'
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('E45').Value = '  +3.27%  '
$ws.Range('D46').Value = "'32.61"
$ws.Range('E46').Value = '  +1.41%  '
$ws.Range('E47').Value = '  +8.63%  '
$ws.Range('D48').Value = "'2.60"
$ws.Range('E48').Value = '  +5.35%  '
$ws.Range('E49').Value = '  +2.29%  '
$ws.Range('D50').Value = "'133.83"
$ws.Range('E50').Value = '  +0.57%  '
$ws.Range('E51').Value = '  -0.05%  '
